$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Doc")

# --- "Max sort:" block: bump the sorting-time / array-size values ---
$ws.Range("C3").Value = 32000
$ws.Range("H3").Value = 32000
$ws.Range("B4").Value = 1000000
$ws.Range("G4").Value = 1000000

# --- "Quick sort:" block ---
# 1) Remove the old "Quick sort:" label row (row 6); everything below shifts up.
$ws.Rows("6").Delete()

# 2) Insert a new blank separator row before the Quick sort data (now starting at row 6).
$ws.Rows("6").Insert()
$ws.Rows("6").Clear()

# 3) Insert a new row for an extra "10,000,000" data point, before the row that now
#    holds the 100,000,000 value (now at row 10).
$ws.Rows("10").Insert()

# --- "Radix sort:" block ---
# 4) Remove the old "Radix sort:" label row (now at row 13).
$ws.Rows("13").Delete()

# 5) Insert a new row for an extra "10,000,000" data point, before the row that now
#    holds the 100,000,000 value (now at row 16).
$ws.Rows("16").Insert()

# --- Fill in the Quick sort block content (rows 7-11) ---
$ws.Range("A7").Value = "Quick sort:"
$ws.Range("F7").Value = "Quick sort:"

$ws.Range("B8").ClearFormats()
$ws.Range("G8").ClearFormats()
$ws.Range("B8").Value = 1000000
$ws.Range("G8").Value = 1000000

$ws.Range("B9").Value = 1000000
$ws.Range("G9").Value = 1000000

$ws.Range("B10").Value = 10000000
$ws.Range("C10").Value = 32000
$ws.Range("G10").Value = 10000000
$ws.Range("H10").Value = 32000

# row 11 (100000000) already holds the right value from the shift; leave as-is.

# --- Fill in the Radix sort block content (rows 13-17) ---
$ws.Range("A13").Value = "Radix sort:"
$ws.Range("F13").Value = "Radix sort:"

$ws.Range("B14").ClearFormats()
$ws.Range("G14").ClearFormats()
$ws.Range("B14").Value = 1000000
$ws.Range("G14").Value = 1000000

$ws.Range("B15").Value = 1000000
$ws.Range("G15").Value = 1000000

$ws.Range("B16").Value = 10000000
$ws.Range("C16").Value = 32000
$ws.Range("G16").Value = 10000000
$ws.Range("H16").Value = 32000

# row 17 (100000000) already holds the right value from the shift; leave as-is.

# --- misc view tweak mirrored from the diff ---
$ws.Range("J15").Select()
